$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New interleaved order of worker statement-of-account rows (periods 1607..2010),
# alternating ROSA ELENA PEDRAZA HENAO and YANETH DUQUE SALAS for each period,
# ascending by period. This replaces the previous layout where each worker's
# periods were listed together in descending order.
$data = @(
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1607",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1607",38640,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1608",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1608",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1609",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1609",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1610",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1610",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1611",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1611",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1612",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1612",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1701",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1701",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1702",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1702",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1703",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1703",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1704",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1704",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1705",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1705",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1706",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1706",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1707",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1707",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1708",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1708",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1709",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1709",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1710",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1710",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1711",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1711",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1712",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1712",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1801",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1801",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1802",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1802",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1803",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1803",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1804",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1804",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1805",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1805",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1806",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1806",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1807",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1807",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1808",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1808",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1809",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1809",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1810",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1810",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1811",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1811",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1812",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1812",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1901",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1901",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1902",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1902",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1903",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1903",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1904",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1904",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1905",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1905",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1906",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1906",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1907",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1907",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1908",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1908",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1909",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1909",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1910",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1910",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1911",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1911",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","1912",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","1912",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2001",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2001",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2002",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2002",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2003",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2003",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2004",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2004",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2005",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2005",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2006",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2006",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2007",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2007",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2008",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2008",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2009",49320,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2009",55200,1380000),
    @("CC","53155281","ROSA ELENA PEDRAZA HENAO","2010",41100,1233000),
    @("CC","45479512","YANETH DUQUE SALAS","2010",46000,1380000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 3).Value = $rowVals[1]
    $ws.Cells.Item($r, 4).Value = $rowVals[2]
    $ws.Cells.Item($r, 5).Value = $rowVals[3]
    $ws.Cells.Item($r, 6).Value = $rowVals[4]
    $ws.Cells.Item($r, 7).Value = $rowVals[5]
}
